$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: Date in A3, time (hours) in B3
$ws.Range("A3").Value = 41549
$ws.Range("A3").NumberFormat = "DD/MM/YY"

$ws.Range("B3").Value = 0.104166666666667
$ws.Range("B3").NumberFormat = "HH:MM:SS"

$ws.Range("B3").Select()
